# color_scheme_palate.pptx — recolor/reposition the existing palette
# rectangles, nudge "Rectangle 8", and add a new small swatch rectangle
# ("Rectangle 1") near the top of the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Recolor the existing palette swatches -----------------------------
$s.Shapes.Item("Rectangle 3").Fill.ForeColor.RGB = 0x7E3A6F   # 46415F -> 6F3A7E
$s.Shapes.Item("Rectangle 4").Fill.ForeColor.RGB = 0x4640A8   # 8D4046 -> A84046
$s.Shapes.Item("Rectangle 5").Fill.ForeColor.RGB = 0x59823A   # 395558 -> 3A8259
$s.Shapes.Item("Rectangle 6").Fill.ForeColor.RGB = 0xA07260   # 576884 -> 6072A0
$s.Shapes.Item("Rectangle 7").Fill.ForeColor.RGB = 0xFFFFFF   # BAB8BD -> FFFFFF
$s.Shapes.Item("Rectangle 9").Fill.ForeColor.RGB = 0x9F0AD4   # 8A5366 -> D40A9F

# "Rectangle 8" also shifts slightly (946425,2668904) -> (9496424,2783204) EMU
$rect8 = $s.Shapes.Item("Rectangle 8")
$rect8.Fill.ForeColor.RGB = 0x59ACC7                          # BE8B5E -> C7AC59
$rect8.Left = 747.7499392598424
$rect8.Top = 219.14992625984252

# --- Add the new swatch rectangle --------------------------------------
# Duplicated from an existing swatch so it inherits the same quick-style
# (lnRef/fillRef/effectRef/fontRef) and text-body shape, then
# repositioned, resized, recolored, cleared of any inherited fill tweaks
# and renamed to match the authored shape.
$newRect = $s.Shapes.Item("Rectangle 3").Duplicate()
$newRect.Name = "Rectangle 1"
$newRect.Left = 179.1
$newRect.Top = 35.100010
$newRect.Width = 56.7
$newRect.Height = 54
$newRect.Fill.ForeColor.RGB = 0xFDF3EA                        # EAF3FD
